# This script updates the cryptocurrency price ("Price", column D) and
# volume-change percentage ("Volume(1h)", column E) figures on the active
# worksheet to reflect the latest scraped values, matching the commit
# "Updated symbol list on Mon Jan 30 09:25:31 UTC 2023 with GitHub Actions".
#
# Each cell holds its number as literal text (e.g. "310.17", "-0.68%"),
# so we explicitly force a text number format before writing the value and
# then restore the cell's original style, ensuring Excel does not silently
# reinterpret the text as a numeric/percentage value and reformat it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "310.17" },
    @{ Cell = "E2"; Value = "-0.68%" },
    @{ Cell = "D3"; Value = "38.41" },
    @{ Cell = "E3"; Value = "-2.18%" },
    @{ Cell = "D4"; Value = "5.079" },
    @{ Cell = "E4"; Value = "-1.14%" },
    @{ Cell = "D5"; Value = "0.07975" },
    @{ Cell = "E5"; Value = "-1.90%" },
    @{ Cell = "D6"; Value = "2.037" },
    @{ Cell = "E6"; Value = "2.28%" },
    @{ Cell = "D7"; Value = "4.459" },
    @{ Cell = "E7"; Value = "5.31%" },
    @{ Cell = "D8"; Value = "8.304" },
    @{ Cell = "E8"; Value = "2.11%" },
    @{ Cell = "D9"; Value = "3.112" },
    @{ Cell = "E9"; Value = "-7.86%" },
    @{ Cell = "D10"; Value = "0.9330" },
    @{ Cell = "E10"; Value = "0.76%" },
    @{ Cell = "D11"; Value = "0.1282" },
    @{ Cell = "E11"; Value = "-8.63%" },
    @{ Cell = "D12"; Value = "0.1904" },
    @{ Cell = "E12"; Value = "-1.25%" },
    @{ Cell = "D13"; Value = "0.08850" },
    @{ Cell = "E13"; Value = "-2.38%" },
    @{ Cell = "D14"; Value = "0.03462" },
    @{ Cell = "E14"; Value = "-1.91%" },
    @{ Cell = "D15"; Value = "0.09693" },
    @{ Cell = "E15"; Value = "-1.22%" },
    @{ Cell = "D16"; Value = "0.001411" },
    @{ Cell = "E16"; Value = "1.01%" },
    @{ Cell = "D17"; Value = "0.006357" },
    @{ Cell = "E17"; Value = "7.23%" },
    @{ Cell = "D18"; Value = "3.581" },
    @{ Cell = "E18"; Value = "-4.73%" },
    @{ Cell = "D19"; Value = "0.3403" },
    @{ Cell = "E19"; Value = "-1.50%" },
    @{ Cell = "D20"; Value = "0.1296" },
    @{ Cell = "E20"; Value = "-1.24%" },
    @{ Cell = "D21"; Value = "5.035" },
    @{ Cell = "E21"; Value = "8.25%" },
    @{ Cell = "D22"; Value = "0.2535" },
    @{ Cell = "E22"; Value = "4.55%" },
    @{ Cell = "D23"; Value = "0.04379" },
    @{ Cell = "E23"; Value = "0.40%" },
    @{ Cell = "D24"; Value = "0.001247" },
    @{ Cell = "D25"; Value = "0.004678" },
    @{ Cell = "E25"; Value = "-2.57%" },
    @{ Cell = "D26"; Value = "0.0003593" },
    @{ Cell = "E26"; Value = "176.36%" },
    @{ Cell = "D39"; Value = "0.02184" },
    @{ Cell = "E39"; Value = "2.53%" },
    @{ Cell = "D40"; Value = "0.05100" },
    @{ Cell = "E40"; Value = "-1.45%" },
    @{ Cell = "D41"; Value = "0.007611" },
    @{ Cell = "E41"; Value = "2.18%" },
    @{ Cell = "D42"; Value = "0.009851" },
    @{ Cell = "E42"; Value = "0.80%" },
    @{ Cell = "D43"; Value = "0.1379" },
    @{ Cell = "E43"; Value = "0.90%" },
    @{ Cell = "D44"; Value = "0.002042" },
    @{ Cell = "E44"; Value = "-4.15%" },
    @{ Cell = "D45"; Value = "0.008839" },
    @{ Cell = "E45"; Value = "-9.41%" },
    @{ Cell = "D46"; Value = "0.00006658" },
    @{ Cell = "E46"; Value = "4.00%" },
    @{ Cell = "D47"; Value = "0.00000000755" },
    @{ Cell = "E47"; Value = "0.61%" },
    @{ Cell = "D48"; Value = "0.003020" },
    @{ Cell = "E48"; Value = "18.17%" },
    @{ Cell = "D49"; Value = "0.001208" },
    @{ Cell = "E49"; Value = "20.74%" },
    @{ Cell = "D50"; Value = "0.00002113" },
    @{ Cell = "E50"; Value = "0.61%" },
    @{ Cell = "D51"; Value = "0.0002013" },
    @{ Cell = "E51"; Value = "0.61%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $originalStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = $originalStyle
}
